$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column D width: 12 -> 13 (Excel's ColumnWidth property adds a constant
# ~0.8333 padding offset versus the raw OOXML <col width> value, so we
# subtract it back out to land exactly on the target width of 13)
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666

# Row 2 (OTROS): VENTA / POR CUMPLIR
$ws.Range("D2").Value = 1460.22
$ws.Range("E2").Value = -1460.22

# Row 3 (PORCELANATO): VENTA / POR CUMPLIR / CUMPLIMIENTO
$ws.Range("D3").Value = 809.04
$ws.Range("E3").Value = 12914.3
$ws.Range("F3").Value = 0.05895357835628935

# Row 4 (TOTAL): VENTA / POR CUMPLIR / CUMPLIMIENTO
$ws.Range("D4").Value = 2269.26
$ws.Range("E4").Value = 11454.08
$ws.Range("F4").Value = 0.165357704465531
